# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Row 16 and Row 17 hold "Periodo Mora" (col E) / "Valor Mora" (col F) data
# that needs to be swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Periodo Mora" values (col E) between row 16 and row 17
$ws.Range("E16").Value = "1906"
$ws.Range("E17").Value = "1907"

# Swap "Valor Mora" values (col F) between row 16 and row 17
$ws.Range("F16").Value = 33600
$ws.Range("F17").Value = 48000
